$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the observation photo path strings to use the new
# "patient-observations/..." prefix instead of "observations/...".
for ($row = 4; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 20)  # column T
    $current = $cell.Value2
    $cell.Value2 = $current -replace '^observations/', 'patient-observations/'
}
